# Sample Project / Main.xlsx - "Rules" sheet
#
# The rule row that used to read "R40" (row 11, column B) is retyped as the
# text "1". The old label is replaced - a brand-new shared string ("1") is
# introduced and cell B11 now points at it - while every other cell, the
# row's style (s="23") and the rest of the sheet stay exactly as they were.
#
# Cells.Item(11, 2) / Range("B11") is the cell in question.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$target = $ws.Range("B11")

# A plain `$target.Value = "1"` would be stored as the *number* 1 (no
# shared-string entry, t="s" would be missing). To reproduce a typed text
# "1" (as in the original edit) without disturbing B11's existing style,
# stage the text in a scratch cell far outside the used range, using
# TEXT() so the staged value is unambiguously text, then copy only the
# *value* (PasteSpecial xlPasteValues = -4163) onto B11. That swaps the
# cell's contents for the text "1" while leaving its formatting/style
# (and every other cell) untouched.
$scratch = $ws.Range("ZZ2")
$scratch.Formula = "=TEXT(1,""0"")"

$scratch.Copy()
$target.PasteSpecial(-4163)

# Clean up the scratch cell and clipboard marquee so no trace of the helper
# cell remains in the saved workbook.
$scratch.Clear()
$excel.CutCopyMode = $false
